# Refresh cryptos list (prices / 1h volume %) per GitHub Actions run.
# Numeric-looking "Price" strings are prefixed with a leading apostrophe so
# Excel stores them as literal text (preserving trailing zeros / multi-dot
# formatting such as "41.50" or "28.164.14") instead of silently coercing
# them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.164.14"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.878.18"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'313.87"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D7").Value = "'0.5133"
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("D8").Value = "'0.3907"
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("D9").Value = "'0.08353"
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("D10").Value = "'1.121"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("D11").Value = "'41.50"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").Value = "'6.227"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "1.884.18"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").Value = "'20.68"
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("D15").Value = "'7.256"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "'1.006"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "'0.00001100"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "'91.11"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "'0.06672"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").Value = "'17.80"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "'6.037"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").Value = "28.205.62"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "'11.13"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("D25").Value = "'2.269"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").Value = "2.094.59"
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.486"
$ws.Range("E27").Value = "  -2.51%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'159.71"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").Value = "'20.64"
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("D30").Value = "'125.10"
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").Value = "'0.1059"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").Value = "'5.854"
$ws.Range("E33").Value = "  +4.64%  "
$ws.Range("D34").Value = "'3.605"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").Value = "'9.662"
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("D36").Value = "'0.02444"
$ws.Range("E36").Value = "  +1.50%  "
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("D38").Value = "'0.2187"
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("D39").Value = "'1.201"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "'0.6497"
$ws.Range("E40").Value = "  +1.56%  "
$ws.Range("D41").Value = "'4.995"
$ws.Range("E41").Value = "  +2.70%  "
$ws.Range("D42").Value = "'1.221"
$ws.Range("E42").Value = "  -1.71%  "
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("D44").Value = "'0.6135"
$ws.Range("E44").Value = "  +0.99%  "
$ws.Range("D45").Value = "'13.14"
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("D46").Value = "'1.283"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").Value = "'3.671"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").Value = "'2.019"
$ws.Range("E48").Value = "  +1.98%  "
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").Value = "'120.81"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").Value = "'78.15"
$ws.Range("E51").Value = "  -0.94%  "
